$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions symbol-list refresh).
# Values are stored as literal text (matching the source inline-string cells),
# so we prefix with an apostrophe to force text entry, then clear the resulting
# quote-prefix formatting so the cell style stays at the sheet default.

$ws.Range("D2").Value = "'321.00"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'-3.45%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'42.59"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'-6.71%"
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'5.217"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'-5.54%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.08228"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-3.57%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'4.317"
$ws.Range("D6").ClearFormats()
$ws.Range("D7").Value = "'1.778"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-13.74%"
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.9494"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-4.19%"
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'-3.09%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1892"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'-1.33%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.09376"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-4.57%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.04625"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-2.00%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'7.451"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-21.40%"
$ws.Range("E13").ClearFormats()
$ws.Range("E14").Value = "'-0.05%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.001299"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-0.37%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'0.005766"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-2.61%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.362"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'-0.64%"
$ws.Range("E17").ClearFormats()
$ws.Range("D19").Value = "'0.3366"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'0.30%"
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'1.09%"
$ws.Range("E20").ClearFormats()
$ws.Range("E21").Value = "'-0.20%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.04155"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'0.22%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.001249"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-4.04%"
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'0.004288"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-6.79%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.0001220"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'-6.29%"
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'0.0002976"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'-0.36%"
$ws.Range("E26").ClearFormats()
$ws.Range("D38").Value = "'0.02681"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'-2.87%"
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.05598"
$ws.Range("D39").ClearFormats()
$ws.Range("D40").Value = "'0.008165"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'3.46%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.1400"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-2.32%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.006545"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-9.76%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.002031"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-4.17%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.007661"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-5.43%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.3480"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-2.17%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006767"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-4.11%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.003075"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'-11.05%"
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'15.87%"
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'-0.26%"
$ws.Range("E51").ClearFormats()
